$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 17618.285
$ws.Range("I11").Value = 17618.285
$ws.Range("K11").Value = 17618.285
$ws.Range("M11").Value = -17478.285
$ws.Range("H28").Value = 597
$ws.Range("J28").Value = 380
$ws.Range("L28").Value = 380
$ws.Range("N28").Value = -1350
$ws.Range("H43").Value = 8335124
$ws.Range("I43").Value = 100000000
$ws.Range("J43").Value = 1953.6364
$ws.Range("K43").Value = 100000000
$ws.Range("L43").Value = 1953.6364
$ws.Range("M43").Value = -99999931
$ws.Range("N43").Value = -2091.6364
$ws.Range("H62").Value = 4990.6665
$ws.Range("I62").Value = 4990.6665
$ws.Range("K62").Value = 4990.6665
$ws.Range("M62").Value = -4366.6665
$ws.Range("H64").Value = 4256.7144
$ws.Range("I64").Value = 4159.6
$ws.Range("K64").Value = 4159.6
$ws.Range("M64").Value = -3911.6
$ws.Range("H65").Value = 4990.6665
$ws.Range("I65").Value = 4990.6665
$ws.Range("K65").Value = 24953.3325
$ws.Range("M65").Value = -21833.3325
$ws.Range("H67").Value = 4256.7144
$ws.Range("I67").Value = 4159.6
$ws.Range("K67").Value = 4159.6
$ws.Range("M67").Value = -3301.6
$ws.Range("H101").Value = 2063.7
$ws.Range("I101").Value = 2140.5
$ws.Range("J101").Value = 1948.5
$ws.Range("K101").Value = 6421.5
$ws.Range("L101").Value = 5845.5
$ws.Range("M101").Value = -4799.5
$ws.Range("N101").Value = -9089.5
$ws.Range("H116").Value = 4641.769
$ws.Range("I116").Value = 5028.9
$ws.Range("K116").Value = 5028.9
$ws.Range("M116").Value = -1586.9
$ws.Range("H135").Value = 199
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 280689.44
$ws.Range("I32").Value = 345372.75
$ws.Range("K32").Value = 345372.75
$ws.Range("M32").Value = -345085.75
$ws.Range("H39").Value = 22500000
$ws.Range("I39").Value = 22500000
$ws.Range("K39").Value = 22500000
$ws.Range("M39").Value = -22499480
$ws.Range("H97").Value = 10062.267
$ws.Range("I97").Value = 11867.5
$ws.Range("K97").Value = 11867.5
$ws.Range("M97").Value = -11371.5
$ws.Range("H132").Value = 2606.2222
$ws.Range("I132").Value = 1570.8889
$ws.Range("J132").Value = 4159.222
$ws.Range("K132").Value = 4712.6667
$ws.Range("L132").Value = 12477.666
$ws.Range("M132").Value = -2182.6667
$ws.Range("N132").Value = -17537.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1242.3529
$ws.Range("J64").Value = 1596.091
$ws.Range("L64").Value = 1596.091
$ws.Range("N64").Value = -2046.091
$ws.Range("H67").Value = 1242.3529
$ws.Range("J67").Value = 1596.091
$ws.Range("L67").Value = 1596.091
$ws.Range("N67").Value = -3156.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11906593
$ws.Range("I16").Value = 20409088
$ws.Range("K16").Value = 20409088
$ws.Range("M16").Value = -20408801
$ws.Range("H58").Value = 2824.6
$ws.Range("I58").Value = 1061.75
$ws.Range("J58").Value = 3999.8333
$ws.Range("K58").Value = 1061.75
$ws.Range("L58").Value = 3999.8333
$ws.Range("M58").Value = -858.75
$ws.Range("N58").Value = -4405.8333
$ws.Range("H107").Value = 2345.25
$ws.Range("I107").Value = 2380.7144
$ws.Range("K107").Value = 2380.7144
$ws.Range("M107").Value = -460.7143999999998
$ws.Range("H113").Value = 11906593
$ws.Range("I113").Value = 20409088
$ws.Range("K113").Value = 20409088
$ws.Range("M113").Value = -20406918
$ws.Range("H134").Value = 2418.5757
$ws.Range("I134").Value = 2435.1667
$ws.Range("J134").Value = 2374.3333
$ws.Range("K134").Value = 7305.500100000001
$ws.Range("L134").Value = 7122.999899999999
$ws.Range("M134").Value = -4770.500100000001
$ws.Range("N134").Value = -12192.9999
$ws.Range("H136").Value = 2824.6
$ws.Range("I136").Value = 1061.75
$ws.Range("J136").Value = 3999.8333
$ws.Range("K136").Value = 3185.25
$ws.Range("L136").Value = 11999.4999
$ws.Range("M136").Value = -635.25
$ws.Range("N136").Value = -17099.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 8500
$ws.Range("J74").Value = 12000
$ws.Range("L74").Value = 36000
$ws.Range("N74").Value = -38122
$ws.Range("H77").Value = 8500
$ws.Range("J77").Value = 12000
$ws.Range("L77").Value = 108000
$ws.Range("N77").Value = -118608
$ws.Range("H122").Value = 6062038.5
$ws.Range("I122").Value = 8333891.5
$ws.Range("K122").Value = 75005023.5
$ws.Range("M122").Value = -75002573.5
$ws.Range("H131").Value = 5350340.5
$ws.Range("J131").Value = 3135.3635
$ws.Range("L131").Value = 9406.0905
$ws.Range("N131").Value = -19486.0905
$ws.Range("H137").Value = 3129
$ws.Range("I137").Value = 2327.1428
$ws.Range("K137").Value = 6981.428400000001
$ws.Range("M137").Value = -1881.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 72.882355
$ws.Range("I2").Value = 71.5
$ws.Range("K2").Value = 71.5
$ws.Range("M2").Value = 41.5
$ws.Range("H38").Value = 15000
$ws.Range("J38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("N38").Value = -15926
$ws.Range("H64").Value = 59999.668
$ws.Range("J64").Value = 59999.668
$ws.Range("L64").Value = 59999.668
$ws.Range("N64").Value = -60495.668
$ws.Range("H67").Value = 59999.668
$ws.Range("J67").Value = 59999.668
$ws.Range("L67").Value = 59999.668
$ws.Range("N67").Value = -61715.668
$ws.Range("H102").Value = 25001664
$ws.Range("I102").Value = 31251422
$ws.Range("J102").Value = 2631
$ws.Range("K102").Value = 31251422
$ws.Range("L102").Value = 2631
$ws.Range("M102").Value = -31249800
$ws.Range("N102").Value = -5875
$ws.Range("H132").Value = 1432746
$ws.Range("I132").Value = 10284
$ws.Range("K132").Value = 30852
$ws.Range("M132").Value = -28322

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6248
$ws.Range("I7").Value = 2166.0667
$ws.Range("K7").Value = 2166.0667
$ws.Range("M7").Value = -2054.0667
$ws.Range("H16").Value = 892.8333
$ws.Range("I16").Value = 838
$ws.Range("J16").Value = 1167
$ws.Range("K16").Value = 838
$ws.Range("L16").Value = 1167
$ws.Range("M16").Value = -668
$ws.Range("N16").Value = -1507
$ws.Range("H22").Value = 4941.7144
$ws.Range("I22").Value = 1998.4
$ws.Range("J22").Value = 5861.5
$ws.Range("K22").Value = 1998.4
$ws.Range("L22").Value = 5861.5
$ws.Range("M22").Value = -1703.4
$ws.Range("N22").Value = -6451.5
$ws.Range("H27").Value = 4941.7144
$ws.Range("I27").Value = 1998.4
$ws.Range("J27").Value = 5861.5
$ws.Range("K27").Value = 1998.4
$ws.Range("L27").Value = 5861.5
$ws.Range("M27").Value = -1891.4
$ws.Range("N27").Value = -6075.5
$ws.Range("H30").Value = 3259.6667
$ws.Range("I30").Value = 3259.6667
$ws.Range("K30").Value = 3259.6667
$ws.Range("M30").Value = -3151.6667
$ws.Range("H93").Value = 2495.5557
$ws.Range("I93").Value = 1714.7333
$ws.Range("J93").Value = 6399.6665
$ws.Range("K93").Value = 1714.7333
$ws.Range("L93").Value = 6399.6665
$ws.Range("M93").Value = -466.7333000000001
$ws.Range("N93").Value = -8895.666499999999
$ws.Range("H126").Value = 6248
$ws.Range("I126").Value = 2166.0667
$ws.Range("K126").Value = 6498.2001
$ws.Range("M126").Value = -4028.2001
$ws.Range("H140").Value = 111666.664
$ws.Range("J140").Value = 111666.664
$ws.Range("L140").Value = 111666.664
$ws.Range("N140").Value = -122026.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 35999.832
$ws.Range("I96").Value = 3999.6667
$ws.Range("J96").Value = 68000
$ws.Range("K96").Value = 3999.6667
$ws.Range("L96").Value = 68000
$ws.Range("M96").Value = -2626.6667
$ws.Range("N96").Value = -70746
$ws.Range("H132").Value = 2445.7036
$ws.Range("J132").Value = 2603.9333
$ws.Range("L132").Value = 7811.7999
$ws.Range("N132").Value = -12871.7999
